# Added Test Data For Hungary/Russia/Finland Market
#
# Adds three new market worksheets (Russia, Finland, Hungary) cloned from the
# existing "Austria"/"Denmark" templates, fills in their market name + user
# story cells, fixes up Finland's repeater list (MZXSDR240 row moved above
# PR1DS/PR8AS), and makes Hungary (the last tab) the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Russia  (clone of "Austria": 19 rows, no MZXSDR240 row)
# ---------------------------------------------------------------------
$austria = $wb.Worksheets.Item("Austria")
$austria.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$russia = $wb.Worksheets.Item($wb.Worksheets.Count)
$russia.Name = "Russia"
# Set the "User Story" cell before the "Description" cell so the new shared
# strings land in the same order as the source edit (ticket id, then market).
$russia.Range("B4").Value = "NGC-2929/T2917"
$russia.Range("B2").Value = "Russia Market"
$russia.Rows.Item(3).RowHeight = 28.8
$russia.Rows.Item(4).RowHeight = 28.8
$russia.Rows.Item(5).RowHeight = 28.8

# ---------------------------------------------------------------------
# Finland  (clone of "Denmark": has the MZXSDR240 row, but it needs to move
# up above PR1DS/PR8AS instead of sitting below them)
# ---------------------------------------------------------------------
$denmark = $wb.Worksheets.Item("Denmark")
$denmark.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$finland = $wb.Worksheets.Item($wb.Worksheets.Count)
$finland.Name = "Finland"
$finland.Range("B4").Value = "NGC-3130/T2887"
$finland.Range("B2").Value = "Finland Market"
$finland.Rows.Item(3).RowHeight = 28.8
$finland.Rows.Item(4).RowHeight = 28.8
$finland.Rows.Item(5).RowHeight = 28.8

# Move "MZXSDR240" from row 18 (after PR1DS/PR8AS) up to a fresh row 16
# (before PR1DS/PR8AS) to match the target layout.
$finland.Cells.Item(18, 1).ClearContents()
$finland.Rows.Item(16).Insert()
$finland.Cells.Item(17, 1).Copy($finland.Cells.Item(16, 1))
$finland.Cells.Item(16, 1).Value = "MZXSDR240"

# ---------------------------------------------------------------------
# Hungary  (clone of "Austria": 19 rows, no MZXSDR240 row) - ends up as the
# last / active tab, same as in the source edit.
# ---------------------------------------------------------------------
$austria2 = $wb.Worksheets.Item("Austria")
$austria2.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"
$hungary.Range("B4").Value = "NGC-3104/T2979"
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Rows.Item(3).RowHeight = 28.8
$hungary.Rows.Item(4).RowHeight = 28.8
$hungary.Rows.Item(5).RowHeight = 28.8

# Hungary becomes the active sheet/tab, as in the target workbook.
$hungary.Activate()
